$wb = $excel.ActiveWorkbook

# --- Sheet: mls ---
$ws = $wb.Worksheets.Item("mls")
$ws.Range("S1").Value = "X17"
$ws.Range("R2").Value = 19
$ws.Range("S2").Value = 22
$ws.Range("S3").Value = 2
$ws.Range("R4").Value = 5
$ws.Range("S4").Value = 3
$ws.Range("R5").Value = 17
$ws.Range("S5").Value = 21
$ws.Range("S6").Value = 25
$ws.Range("R7").Value = 20
$ws.Range("S7").Value = 23
$ws.Range("S8").Value = 16
$ws.Range("S9").Value = 28
$ws.Range("R10").Value = 12
$ws.Range("S10").Value = 11
$ws.Range("R11").Value = 6
$ws.Range("S11").Value = 8
$ws.Range("R12").Value = 15
$ws.Range("S12").Value = 17
$ws.Range("R13").Value = 16
$ws.Range("S13").Value = 18
$ws.Range("S14").Value = 1
$ws.Range("R15").Value = 11
$ws.Range("S15").Value = 12
$ws.Range("R16").Value = 23
$ws.Range("S16").Value = 19
$ws.Range("R17").Value = 7
$ws.Range("S17").Value = 9
$ws.Range("R18").Value = 13
$ws.Range("S18").Value = 13
$ws.Range("S19").Value = 7
$ws.Range("R20").Value = 8
$ws.Range("S20").Value = 4
$ws.Range("R21").Value = 10
$ws.Range("S21").Value = 10
$ws.Range("R22").Value = 9
$ws.Range("S22").Value = 5
$ws.Range("R23").Value = 21
$ws.Range("S23").Value = 15
$ws.Range("S24").Value = 6
$ws.Range("R25").Value = 26
$ws.Range("S25").Value = 27
$ws.Range("R26").Value = 14
$ws.Range("S26").Value = 14
$ws.Range("R27").Value = 25
$ws.Range("S27").Value = 26
$ws.Range("R28").Value = 24
$ws.Range("S28").Value = 24
$ws.Range("R29").Value = 18
$ws.Range("S29").Value = 20

# --- Sheet: nor ---
$ws = $wb.Worksheets.Item("nor")
$ws.Range("N1").Value = "X12"
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = 7
$ws.Range("M3").Value = 4
$ws.Range("N3").Value = 5
$ws.Range("N4").Value = 9
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = 13
$ws.Range("N6").Value = 15
$ws.Range("N7").Value = 16
$ws.Range("N8").Value = 1
$ws.Range("N9").Value = 2
$ws.Range("M10").Value = 13
$ws.Range("N10").Value = 10
$ws.Range("M11").Value = 5
$ws.Range("N11").Value = 6
$ws.Range("M12").Value = 10
$ws.Range("N12").Value = 11
$ws.Range("M13").Value = 8
$ws.Range("N13").Value = 8
$ws.Range("M14").Value = 6
$ws.Range("N14").Value = 4
$ws.Range("M15").Value = 11
$ws.Range("N15").Value = 12
$ws.Range("M16").Value = 14
$ws.Range("N16").Value = 14
$ws.Range("N17").Value = 3

# --- Sheet: swe ---
$ws = $wb.Worksheets.Item("swe")
$ws.Range("N1").Value = "X12"
$ws.Range("N2").Value = 1
$ws.Range("N3").Value = 14
$ws.Range("M4").Value = 3
$ws.Range("N4").Value = 3
$ws.Range("M5").Value = 6
$ws.Range("N5").Value = 7
$ws.Range("M6").Value = 10
$ws.Range("N6").Value = 8
$ws.Range("N7").Value = 2
$ws.Range("M8").Value = 4
$ws.Range("N8").Value = 4
$ws.Range("N9").Value = 16
$ws.Range("M10").Value = 5
$ws.Range("N10").Value = 6
$ws.Range("M11").Value = 7
$ws.Range("N11").Value = 5
$ws.Range("N12").Value = 9
$ws.Range("M13").Value = 9
$ws.Range("N13").Value = 10
$ws.Range("N14").Value = 11
$ws.Range("N15").Value = 15
$ws.Range("N16").Value = 13
$ws.Range("N17").Value = 12

# --- Sheet: bra ---
$ws = $wb.Worksheets.Item("bra")
$ws.Range("M1").Value = "X11"
$ws.Range("N1").Value = "X12"
$ws.Range("O1").Value = "X13"
$ws.Range("P1").Value = "X14"
$ws.Range("M2").Value = 8
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 14
$ws.Range("P2").Value = 16
$ws.Range("M3").Value = 5
$ws.Range("N3").Value = 4
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 3
$ws.Range("M4").Value = 6
$ws.Range("N4").Value = 5
$ws.Range("O4").Value = 4
$ws.Range("P4").Value = 4
$ws.Range("M5").Value = 15
$ws.Range("N5").Value = 17
$ws.Range("O5").Value = 12
$ws.Range("P5").Value = 13
$ws.Range("M6").Value = 9
$ws.Range("N6").Value = 7
$ws.Range("O6").Value = 11
$ws.Range("P6").Value = 9
$ws.Range("M7").Value = 16
$ws.Range("N7").Value = 11
$ws.Range("O7").Value = 6
$ws.Range("P7").Value = 10
$ws.Range("M8").Value = 10
$ws.Range("N8").Value = 8
$ws.Range("O8").Value = 7
$ws.Range("P8").Value = 11
$ws.Range("M9").Value = 11
$ws.Range("N9").Value = 12
$ws.Range("O9").Value = 13
$ws.Range("P9").Value = 14
$ws.Range("M10").Value = 2
$ws.Range("N10").Value = 2
$ws.Range("O10").Value = 2
$ws.Range("P10").Value = 2
$ws.Range("M11").Value = 7
$ws.Range("N11").Value = 13
$ws.Range("O11").Value = 15
$ws.Range("P11").Value = 17
$ws.Range("M12").Value = 17
$ws.Range("N12").Value = 18
$ws.Range("O12").Value = 18
$ws.Range("P12").Value = 18
$ws.Range("M13").Value = 18
$ws.Range("N13").Value = 14
$ws.Range("O13").Value = 16
$ws.Range("P13").Value = 12
$ws.Range("M14").Value = 12
$ws.Range("N14").Value = 15
$ws.Range("O14").Value = 8
$ws.Range("P14").Value = 6
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = 20
$ws.Range("O15").Value = 19
$ws.Range("P15").Value = 20
$ws.Range("M16").Value = 13
$ws.Range("N16").Value = 16
$ws.Range("O16").Value = 17
$ws.Range("P16").Value = 15
$ws.Range("M17").Value = 3
$ws.Range("N17").Value = 3
$ws.Range("O17").Value = 5
$ws.Range("P17").Value = 5
$ws.Range("M18").Value = 19
$ws.Range("N18").Value = 19
$ws.Range("O18").Value = 20
$ws.Range("P18").Value = 19
$ws.Range("M19").Value = 1
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 1
$ws.Range("P19").Value = 1
$ws.Range("M20").Value = 14
$ws.Range("N20").Value = 9
$ws.Range("O20").Value = 9
$ws.Range("P20").Value = 7
$ws.Range("M21").Value = 4
$ws.Range("N21").Value = 6
$ws.Range("O21").Value = 10
$ws.Range("P21").Value = 8

# --- Sheet: irl ---
$ws = $wb.Worksheets.Item("irl")
$ws.Range("W1").Value = "X21"
$ws.Range("W2").Value = 6
$ws.Range("W3").Value = 3
$ws.Range("W4").Value = 8
$ws.Range("W5").Value = 2
$ws.Range("W6").Value = 9
$ws.Range("W7").Value = 1
$ws.Range("W8").Value = 7
$ws.Range("W9").Value = 5
$ws.Range("W10").Value = 4
$ws.Range("W11").Value = 10

# --- Sheet: jpn ---
$ws = $wb.Worksheets.Item("jpn")
$ws.Range("T1").Value = "X18"
$ws.Range("T2").Value = 14
$ws.Range("T3").Value = 6
$ws.Range("T4").Value = 8
$ws.Range("T5").Value = 12
$ws.Range("T6").Value = 9
$ws.Range("T7").Value = 15
$ws.Range("T8").Value = 2
$ws.Range("T9").Value = 4
$ws.Range("T10").Value = 3
$ws.Range("T11").Value = 13
$ws.Range("T12").Value = 10
$ws.Range("T13").Value = 7
$ws.Range("T14").Value = 5
$ws.Range("T15").Value = 17
$ws.Range("T16").Value = 16
$ws.Range("T17").Value = 11
$ws.Range("T18").Value = 18
$ws.Range("T19").Value = 1
